$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REVISION")

# Insert a new row above the current row 2 (existing data row shifts down to row 3)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the purchase data extracted from the PDF
$ws.Range("C2").Value = "MX07"

# E2 holds an IMEI that looks numeric but must stay a text value
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "867501041079210"
$ws.Range("E2").ClearFormats()

$ws.Range("F2").Value = 70011903
$ws.Range("G2").Value = "HUAWEI LTE ROUTER B310S-518 BLANCO de Pedro"
$ws.Range("L2").Value = 9512434283
$ws.Range("M2").Value = "VICSA"
$ws.Range("P2").Value = 6363.48
$ws.Range("Q2").Value = "RORP880418621"
$ws.Range("R2").Value = "RODRIGUEZ RODRIGUEZ PEDRO"
$ws.Range("T2").Value = "AV JUVENTUD"
$ws.Range("U2").Value = "S/N"
$ws.Range("V2").Value = "SN"
$ws.Range("W2").Value = "BARR LA PEÑA"
$ws.Range("X2").Value = "ASUNCION NOCHI"
$ws.Range("Y2").Value = 69600
$ws.Range("Z2").Value = "ASUNCION NOCHIXTLAN"
$ws.Range("AA2").Value = "OAXACA"
$ws.Range("AB2").Value = "skaniahome@gmail.com"
$ws.Range("AC2").Value = "M47"
$ws.Range("AD2").Value = 24
$ws.Range("AF2").Value = 2089
$ws.Range("AK2").Value = "PUE"

# AL2 holds a leading-zero code that must stay a text value
$ws.Range("AL2").NumberFormat = "@"
$ws.Range("AL2").Value = "01"
$ws.Range("AL2").ClearFormats()

$ws.Range("AM2").Value = "P01"
